# Nerul Quotation workbook update:
#  - Fill in the previously-blank expense rows 11-20 on Sheet3 (dates, taken-by,
#    item description, amount), reusing the existing row-2 cell formatting.
#  - Extend the sheet with 30 more blank numbered rows (21-50) so the running
#    Sr. No. column keeps going.
#  - Move the active tab/selection over to Sheet3 (it was Sheet4 before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Clone the still-plain formatting of row 20 onto the new blank rows 21-50
# FIRST, while it is still untouched -- otherwise it would pick up the
# restyled look applied to rows 11-20 just below.
$ws.Range("A20:E20").Copy()
$ws.Range("A21:E50").PasteSpecial(-4122)

# Clone the formatting (date style on B, named-person style on C/D, plain
# style on A/E) from the first data row down across the new rows, same as
# Excel's "fill down formatting" behaviour, before any values are entered.
$ws.Range("A2:E2").Copy()
$ws.Range("A11:E20").PasteSpecial(-4122)

# Sr. No. column for the freshly added blank rows.
for ($i = 21; $i -le 50; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Date / Taken by / Amount for the newly populated rows.
$ws.Range("B11").Value = 45183
$ws.Range("C11").Value = "Sandesh"
$ws.Range("E11").Value = 300

$ws.Range("B12").Value = 45197
$ws.Range("C12").Value = "Sandesh"
$ws.Range("E12").Value = 200

$ws.Range("B13").Value = 45197
$ws.Range("C13").Value = "Sandesh"
$ws.Range("E13").Value = 262

$ws.Range("B14").Value = 45172
$ws.Range("C14").Value = "Sandesh"
$ws.Range("E14").Value = 50

$ws.Range("B15").Value = 45173
$ws.Range("C15").Value = "Deepak"
$ws.Range("E15").Value = 728

$ws.Range("B16").Value = 45173
$ws.Range("C16").Value = "Deepak"
$ws.Range("E16").Value = 590

$ws.Range("B17").Value = 45173
$ws.Range("C17").Value = "Deepak"
$ws.Range("E17").Value = 250

$ws.Range("B18").Value = 45173
$ws.Range("C18").Value = "Deepak"
$ws.Range("E18").Value = 500

$ws.Range("B19").Value = 45173
$ws.Range("C19").Value = "Deepak"
$ws.Range("E19").Value = 400

$ws.Range("B20").Value = 45174
$ws.Range("C20").Value = "Sandesh"
$ws.Range("E20").Value = 200

# Item Description (column D) -- entered in this exact order to match the
# original authoring order of the underlying shared-string table.
$ws.Range("D11").Value = "Nerul Material"
$ws.Range("D12").Value = "Petrol (147810)"
$ws.Range("D13").Value = "Nerul Work"
$ws.Range("D19").Value = "Petrol (No Bill)"
$ws.Range("D14").Value = "RJ 45 Connectors"
$ws.Range("D15").Value = "Krishna Paints and Sanitary"
$ws.Range("D16").Value = "Avon Computers"
$ws.Range("D17").Value = "Leth Charges"
$ws.Range("D18").Value = "Azrenkar Hardware"
$ws.Range("D20").Value = "Petrol (85346)"

# Sheet3 becomes the active tab/sheet, with B21 selected (the next empty
# item row) -- previously Sheet4 was the active tab.
$ws.Activate()
$ws.Range("B21").Select()
